# Set forecast confidence default to High, update to v1.176
#
# Recolors a handful of reviewer-comment paragraphs in the "Optimisation"
# section to red (RGB EE0000 -> wdColor value 238), matching the rest of
# the document's red review-comment styling.
#
# wdColor values are packed as 0x00BBGGRR, so a pure-red RGB value of
# EE0000 (R=0xEE, G=0x00, B=0x00) is simply 0x0000EE = 238 decimal.

$d = $word.ActiveDocument
$wdRed = 238

function Colorize-ParagraphContaining($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1)
        $para.Range.Font.Color = $wdRed
    }
}

# "The allocation % ..." / "Default is cash 100% ..." paragraph
Colorize-ParagraphContaining("The allocation % and constrains are for the current holdings.")

# "Show the before and after allocations ..." / "eg" / "100% cash => ..." paragraph
Colorize-ParagraphContaining("Show the before and after allocations as an output")

# "The output isn't correct ... methodology used" paragraph
Colorize-ParagraphContaining("The output isn")

# "Resampled approach is the preference ... very conservative investor." paragraph
Colorize-ParagraphContaining("Resampled approach is the preference")
